# "updated existing failed testcases after launchpad release"
#
# The "Quiz Dashboard Classes" expected-value rows (F29:F34, method
# verifyQuizDashboardClasses) were missing the newly introduced
# "Class 8-D" and "Class 11-D" sections from their comma-separated class
# list. Bring them in line with the rest of the class roster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "Class 3-A, Class 3-B, Class 3-C, Class 4-A, Class 4-B, Class 4-C, Class 5-A, Class 5-B, Class 5-C, Class 6-A, Class 6-B, Class 6-C, Class 7-A, Class 7-B, Class 7-C, Class 8-A, Class 8-B, Class 8-C, Class 8-D, Class 9-A, Class 9-B, Class 9-C, Class 10-A, Class 10-B, Class 10-C, Class 11-A, Class 11-B, Class 11-C, Class 11-D, Class 12-A, Class 12-B, Class 12-C"

$ws.Range("F29:F34").Value = $newValue

# Reflect where the edit happened in the saved view state.
$ws.Range("F29:F34").Select()
$excel.ActiveWindow.ScrollRow = 24
